# Update recalculated TPM-derived metrics in the NATMI LR-pairs output sheet.
# Values below reflect a re-run of the scoring pipeline with updated TPM input
# (see commit message: "update scripts wuth new tpm").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Sending=FAPs, Ligand=Lgi3, Receptor=Stx1a, Target=ECs)
$ws.Range("G2").Value = 0.2994946666666667
$ws.Range("H2").Value = 0.8984840000000001
$ws.Range("J2").Value = 0.4989451716962828
$ws.Range("M2").Value = 1.542131666666667
$ws.Range("N2").Value = 4.626395
$ws.Range("O2").Value = 0.2639921135523384
$ws.Range("P2").Value = 0.2639921135523383
$ws.Range("Q2").Value = 0.4618602094644446
$ws.Range("R2").Value = 4.156741885180001
$ws.Range("S2").Value = 0.131717590422836
$ws.Range("T2").Value = 0.131717590422836

# Row 3 (Sending=FAPs, Ligand=Lgi3, Receptor=Stx1a, Target=FAPs)
$ws.Range("G3").Value = 0.2994946666666667
$ws.Range("H3").Value = 0.8984840000000001
$ws.Range("J3").Value = 0.4989451716962828
$ws.Range("N3").Value = 7.699008000000001
$ws.Range("O3").Value = 0.4393220626808479
$ws.Range("P3").Value = 0.4393220626808478
$ws.Range("Q3").Value = 0.7686039448746668
$ws.Range("R3").Value = 6.917435503872001
$ws.Range("S3").Value = 0.2191976219942607
$ws.Range("T3").Value = 0.2191976219942607

# Row 4 (Sending=FAPs, Ligand=Lgi3, Receptor=Stx1a, Target=MuSCs)
$ws.Range("G4").Value = 0.2994946666666667
$ws.Range("H4").Value = 0.8984840000000001
$ws.Range("J4").Value = 0.4989451716962828
$ws.Range("O4").Value = 0.2966858237668138
$ws.Range("P4").Value = 0.2966858237668138
$ws.Range("Q4").Value = 0.5190585993884445
$ws.Range("R4").Value = 4.671527394496001
$ws.Range("S4").Value = 0.148029959279186
$ws.Range("T4").Value = 0.148029959279186

# Row 5 (Sending=MuSCs, Ligand=Lgi3, Receptor=Stx1a, Target=ECs)
$ws.Range("I5").Value = 0.5010548283037172
$ws.Range("M5").Value = 1.542131666666667
$ws.Range("N5").Value = 4.626395
$ws.Range("O5").Value = 0.2639921135523384
$ws.Range("P5").Value = 0.2639921135523383
$ws.Range("Q5").Value = 0.4638130621983334
$ws.Range("R5").Value = 4.174317559785
$ws.Range("S5").Value = 0.1322745231295023
$ws.Range("T5").Value = 0.1322745231295023

# Row 6 (Sending=MuSCs, Ligand=Lgi3, Receptor=Stx1a, Target=FAPs)
$ws.Range("I6").Value = 0.5010548283037172
$ws.Range("N6").Value = 7.699008000000001
$ws.Range("O6").Value = 0.4393220626808479
$ws.Range("P6").Value = 0.4393220626808478
$ws.Range("Q6").Value = 0.7718537816960001
$ws.Range("R6").Value = 6.946684035264
$ws.Range("S6").Value = 0.2201244406865871
$ws.Range("T6").Value = 0.2201244406865871

# Row 7 (Sending=MuSCs, Ligand=Lgi3, Receptor=Stx1a, Target=MuSCs)
$ws.Range("I7").Value = 0.5010548283037172
$ws.Range("O7").Value = 0.2966858237668138
$ws.Range("P7").Value = 0.2966858237668138
$ws.Range("S7").Value = 0.1486558644876278
$ws.Range("T7").Value = 0.1486558644876278
